$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for numeric-looking price values so Excel does not
# auto-convert them to numbers (which would lose exact text such as trailing zeros).
$textCells = @("D5", "D6", "D11", "D15", "D17", "D19", "D21", "D22", "D23", "D27", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D41", "D42", "D45", "D47", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "41.690.64"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.478.53"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "319.31"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "92.66"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D11").Value = "33.18"
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "2.861.15"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "15.57"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "2.466.30"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "0.797"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "41.642.07"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "70.85"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "11.25"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "239.70"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "25.02"
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "36.71"
$ws.Range("E30").Value = "  +5.23%  "
$ws.Range("D31").Value = "157.47"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").Value = "5.44"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "0.0766"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "17.18"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").Value = "  +3.46%  "
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "1.997.28"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "18.77"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("D47").Value = "9.53"
$ws.Range("E47").Value = "  +6.97%  "
$ws.Range("D48").Value = "2.719.82"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "98.08"
$ws.Range("D50").Value = "75.58"
$ws.Range("E50").Value = "  +5.60%  "
$ws.Range("D51").Value = "67.15"
$ws.Range("E51").Value = "  +1.33%  "
